$wb = $excel.ActiveWorkbook

$small = $wb.Worksheets.Item("SmallRobotPerimeter")

# SmallRobotPerimeter is no longer the active tab; its lingering selection
# becomes the whole A1:A10 block instead of the old A7 row selection.
$small.Range("A1:A10").Select() | Out-Null

# Add the new sheet right after SmallRobotPerimeter
$ws = $wb.Worksheets.Add($null, $small)
$ws.Name = "BigRobotPerimeter"

$ws.Range("A1").Value = "Bottom Cut"
$ws.Range("A2").Value = 120
$ws.Range("A3").Value = 127
$ws.Range("A4").Value = 130
$ws.Range("A5").Value = 300
$ws.Range("A6").Formula = "=A4"
$ws.Range("A7").Formula = "=A3"
$ws.Range("A10").Formula = "=SUM(A2:A7)"

$ws.Columns("A").ColumnWidth = 30.28515625

$ws.Range("A6").Select() | Out-Null
